$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 201
$ws1.Range("F4").Value = 130

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 201
$ws4.Range("F4").Value = 130
